# Insert a new data row at row 624 on Sheet1, shifting existing rows
# 624-731 down to 625-732, and populate the new row with the new
# price-record data (date 2023-10-19, calidad "Primera", origen "Perú", etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 624 (pushes existing row 624.. down by one)
$ws.Rows.Item(624).Insert()

# Populate the new row 624 with the new record's values
$ws.Cells.Item(624, 1).Value = 3
$ws.Cells.Item(624, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(624, 3).Value = "Coquimbo"
$ws.Cells.Item(624, 4).Value = 45218
$ws.Cells.Item(624, 5).Value = 5
$ws.Cells.Item(624, 6).Value = 100112028
$ws.Cells.Item(624, 7).Value = "Sandia"
$ws.Cells.Item(624, 8).Value = "Sin especificar"
$ws.Cells.Item(624, 9).Value = "Primera"
$ws.Cells.Item(624, 10).Value = 2800
$ws.Cells.Item(624, 11).Value = 650
$ws.Cells.Item(624, 12).Value = 700
$ws.Cells.Item(624, 13).Value = 677
$ws.Cells.Item(624, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(624, 15).Value = "Perú"
$ws.Cells.Item(624, 16).Value = 677
$ws.Cells.Item(624, 17).Value = 1
$ws.Cells.Item(624, 18).Value = "Hortaliza"
